$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '65.731.12'
$ws.Range('E2').Value = '  +0.85%  '

$ws.Range('D3').Value = '3.403.33'
$ws.Range('E3').Value = '  -0.28%  '

$ws.Range('E4').Value = '  -0.14%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '561.96'
$ws.Range('E5').Value = '  -0.19%  '

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '177.01'
$ws.Range('E6').Value = '  -0.03%  '

$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.632'
$ws.Range('E7').Value = '  +0.36%  '

$ws.Range('D8').Value = '3.393.14'
$ws.Range('E8').Value = '  -0.18%  '

$ws.Range('E9').Value = '  -0.06%  '

$ws.Range('E10').Value = '  +4.67%  '

$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.638'
$ws.Range('E11').Value = '  +0.25%  '

$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '53.78'
$ws.Range('E12').Value = '  -2.01%  '

$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.0000279'
$ws.Range('E13').Value = '  +0.22%  '

$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '9.22'
$ws.Range('E14').Value = '  +0.36%  '

$ws.Range('D15').Value = '3.936.95'
$ws.Range('E15').Value = '  -0.55%  '

$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '18.35'
$ws.Range('E16').Value = '  -0.18%  '

$ws.Range('B17').Value = 'WrappedEther'
$ws.Range('C17').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D17').Value = '3.415.31'
$ws.Range('E17').Value = '  -0.32%  '

$ws.Range('B18').Value = 'TRON'
$ws.Range('C18').Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.119'
$ws.Range('E18').Value = '  +0.58%  '

$ws.Range('D19').Value = '65.704.49'
$ws.Range('E19').Value = '  +1.29%  '

$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '11.87'
$ws.Range('E20').Value = '  -0.68%  '

$ws.Range('E21').Value = '  +0.69%  '

$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '485.02'
$ws.Range('E22').Value = '  +2.68%  '

$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '4.94'
$ws.Range('E23').Value = '  -1.40%  '

$ws.Range('B24').Value = 'Litecoin'
$ws.Range('C24').Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '90.14'
$ws.Range('E24').Value = '  +3.59%  '

$ws.Range('B25').Value = 'PancakeSwap'
$ws.Range('C25').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '4.12'
$ws.Range('E25').Value = '  -0.94%  '

$ws.Range('B26').Value = 'InternetComputer(DFINITY)'
$ws.Range('C26').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '14.29'
$ws.Range('E26').Value = '  +4.06%  '

$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '2.93'
$ws.Range('E27').Value = '  +1.62%  '

$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '10.69'
$ws.Range('E28').Value = '  -2.16%  '

$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '8.76'
$ws.Range('E29').Value = '  -1.63%  '

$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '31.47'
$ws.Range('E30').Value = '  +2.49%  '

$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '6.57'
$ws.Range('E31').Value = '  -2.13%  '

$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '11.55'
$ws.Range('E32').Value = '  -0.41%  '

$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '62.64'
$ws.Range('E33').Value = '  +3.83%  '

$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '576.98'
$ws.Range('E34').Value = '  -1.27%  '

$ws.Range('E35').Value = '  -1.09%  '

$ws.Range('E36').Value = '  +0.00%  '

$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '3.61'
$ws.Range('E37').Value = '  +3.67%  '

$ws.Range('E38').Value = '  -0.37%  '

$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '36.02'
$ws.Range('E39').Value = '  -0.23%  '

$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.375'
$ws.Range('E40').Value = '  -0.04%  '

$ws.Range('D41').Value = '0.0₃0742'
$ws.Range('E41').Value = '  -2.77%  '

$ws.Range('D42').Value = '3.106.79'
$ws.Range('E42').Value = '  -0.37%  '

$ws.Range('E43').Value = '  -2.11%  '

$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.0419'
$ws.Range('E44').Value = '  +0.91%  '

$ws.Range('E45').Value = '  +0.50%  '

$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '3.18'
$ws.Range('E46').Value = '  -1.12%  '

$ws.Range('E47').Value = '  -3.75%  '

$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.998'
$ws.Range('E48').Value = '  -0.28%  '

$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '140.11'
$ws.Range('E49').Value = '  +2.57%  '

$ws.Range('E50').Value = '  +0.14%  '

$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '8.48'
$ws.Range('E51').Value = '  -0.41%  '
